{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document currently starts with the \"Write Up\" title paragraph,\n// immediately followed by an (empty) paragraph. Insert the three new\n// paragraphs right before that first empty paragraph, so they land right\n// after the title, in document order.\nconst anchorPara = paragraphs.items[1];\n\nconst para1 = anchorPara.insertParagraph(\n  \"This week, we will be starting a whole new series. This time we will be looking into JavaScript with our Visual Studio Code compiler, and manipulating the DOM. The DOM stands for the document object module, and we can use JavaScript to make changes to it as opposed to making our changes in HTML.\",\n  Word.InsertLocation.before\n);\n\nconst para2 = anchorPara.insertParagraph(\n  \"So, if this is something that you might be interested in learning just a bit more about, then please join us for our brand-new article this week entitled:\",\n  Word.InsertLocation.before\n);\n\nconst para3 = anchorPara.insertParagraph(\n  \"1 Append Text to the Body\",\n  Word.InsertLocation.before\n);\npara3.style = \"Heading1\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document currently starts with the \"Write Up\" title paragraph,\n# immediately followed by an (empty) paragraph. Insert the three new\n# paragraphs right before that first empty paragraph, so they land right\n# after the title, in document order.\n$anchor = $d.Paragraphs.Item(2)\n$anchor.Range.InsertParagraphBefore()\n\n$p1 = $d.Paragraphs.Item(2)\n$p1.Range.InsertAfter(\"This week, we will be starting a whole new series. This time we will be looking into JavaScript with our Visual Studio Code compiler, and manipulating the DOM. The DOM stands for the document object module, and we can use JavaScript to make changes to it as opposed to making our changes in HTML.\")\n\n$anchor2 = $d.Paragraphs.Item(3)\n$anchor2.Range.InsertParagraphBefore()\n\n$p2 = $d.Paragraphs.Item(3)\n$p2.Range.InsertAfter(\"So, if this is something that you might be interested in learning just a bit more about, then please join us for our brand-new article this week entitled:\")\n\n$anchor3 = $d.Paragraphs.Item(4)\n$anchor3.Range.InsertParagraphBefore()\n\n$p3 = $d.Paragraphs.Item(4)\n$p3.Range.InsertAfter(\"1 Append Text to the Body\")\n$p3.Style = \"Heading1\"\n"}
